$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 previously held the "Aspirator pentru piscina 1,2 m" product; replace it
# with the "Cafea Organica" product data that used to live in row 6.
$ws.Range("A2").Value = "Cafea Organica House Roast, Exhale, boabe"
$ws.Range("B2").Value = "https://www.emag.ro/cafea-organica-house-roast-exhale-boabe-x001qp1mud/pd/DRLBHKYBM"

# Row 6 is now redundant (its data moved up to row 2), so clear it out entirely.
$ws.Range("A6:B6").ClearContents()

# Update the active selection like Excel would after this edit.
$ws.Range("A6").Select()
